# Trade #17 closed at 2026-02-17 20:04:03 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1399.78
$wsSummary.Range("B4").Value = -0.23
$wsSummary.Range("B5").Value = -0.27
$wsSummary.Range("B6").Value = 17
$wsSummary.Range("B7").Value = 8
$wsSummary.Range("B9").Value = 47.06

# --- Strategy Status sheet ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C5").Value = 99.78
$wsStatus.Range("D5").Value = 17
$wsStatus.Range("E5").Value = -0.23
$wsStatus.Range("F5").Value = -0.22
$wsStatus.Range("G5").Value = 47.06

# --- New trade row (Trade #17) shared by "All Trades" and "MarketMaking" sheets ---
$newRow = @{
    A = 17
    B = "2026-02-17"
    C = "20:03:56"
    D = "MarketMaking"
    E = "UP"
    F = 0.9399999999999999
    G = 0.95
    H = "CLOSED"
    I = 1.0638
    J = 0.01
    K = 99.78
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.1
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A18").Value = $newRow.A

    # Force the date-looking string to stay plain text (not auto-parsed into a
    # date serial) without leaving a quote-prefixed / text-formatted style on
    # the cell: apply a text format just long enough to take the literal
    # value, then snap the style back to the workbook default.
    $ws.Range("B18").NumberFormat = "@"
    $ws.Range("B18").Value = $newRow.B
    $ws.Range("B18").Style = "Normal"

    $ws.Range("C18").Value = $newRow.C
    $ws.Range("D18").Value = $newRow.D
    $ws.Range("E18").Value = $newRow.E
    $ws.Range("F18").Value = $newRow.F
    $ws.Range("G18").Value = $newRow.G
    $ws.Range("H18").Value = $newRow.H
    $ws.Range("I18").Value = $newRow.I
    $ws.Range("J18").Value = $newRow.J
    $ws.Range("K18").Value = $newRow.K
    $ws.Range("L18").Value = $newRow.L
    $ws.Range("M18").Value = $newRow.M
    $ws.Range("N18").Value = $newRow.N
    $ws.Range("O18").Value = $newRow.O
    $ws.Range("P18").Value = $newRow.P
    $ws.Range("Q18").Value = $newRow.Q
}
